$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D28").Value = '[name="촌장"]  그렇습니다. 제가 바로 벨로니 마을 촌장입니다. 여러분께 길 안내를 해드리겠습니다. 여기에 서 계시지 말고, 이쪽으로 오시죠.
'
$ws.Range("D65").Value = '[name="수르트"]  ……(무언가를 기록하고 있다.)
'
$ws.Range("D70").Value = '[name="메테오라이트"]  ………
'
$ws.Range("D72").Value = '[name="메테오라이트"]  확실히 본 적 있기는 하지만, 단지……
'
$ws.Range("D86").Value = '[name="수르트"]  ……그럴지도.
'
$ws.Range("D102").Value = '[name="메테오라이트"]  하지만, 무언가를 찾기 위해 보이지 않는 앞을 향해 나아갈 땐, 주변의 사물에 대해서도 주의를 기울여야 해.
'
$ws.Range("D144").Value = '[name="수르트"]  이건……
'
$ws.Range("D161").Value = '바스락 바스락… 원석충의 소리다.
'
$ws.Range("D195").Value = '[name="메테오라이트"]  어…… 그 말이 틀린 건 아니긴 한데……
'
$ws.Range("D205").Value = '[name="수르트"]  하지만……
'
